$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "30.689.70"
Set-TextValue "E2" "  +2.09%  "

Set-TextValue "D3" "1.891.64"
Set-TextValue "E3" "  +0.80%  "

Set-TextValue "E4" "  +0.20%  "

Set-TextValue "D5" "245.22"
Set-TextValue "E5" "  +0.70%  "

Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  +0.13%  "

Set-TextValue "D7" "0.4934"
Set-TextValue "E7" "  +0.38%  "

Set-TextValue "D8" "0.2959"
Set-TextValue "E8" "  +1.16%  "

Set-TextValue "D9" "0.06793"
Set-TextValue "E9" "  +2.70%  "

Set-TextValue "D10" "1.887.49"
Set-TextValue "E10" "  +0.52%  "

Set-TextValue "D11" "17.23"
Set-TextValue "E11" "  +4.15%  "

Set-TextValue "D12" "0.07253"
Set-TextValue "E12" "  +0.55%  "

Set-TextValue "D13" "90.89"
Set-TextValue "E13" "  +5.27%  "

Set-TextValue "D14" "0.6788"
Set-TextValue "E14" "  +1.77%  "

Set-TextValue "D15" "5.050"
Set-TextValue "E15" "  +2.67%  "

Set-TextValue "D16" "30.667.90"
Set-TextValue "E16" "  +2.10%  "

Set-TextValue "D17" "0.000007998"
Set-TextValue "E17" "  +2.41%  "

Set-TextValue "E18" "  +0.14%  "

Set-TextValue "E19" "  +2.71%  "

Set-TextValue "D20" "2.131.54"
Set-TextValue "E20" "  +0.37%  "

Set-TextValue "E21" "  +0.40%  "

Set-TextValue "D22" "4.819"
Set-TextValue "E22" "  +0.58%  "

Set-TextValue "D23" "189.72"
Set-TextValue "E23" "  +32.77%  "

Set-TextValue "D24" "6.142"
Set-TextValue "E24" "  +4.75%  "

Set-TextValue "D25" "9.411"
Set-TextValue "E25" "  +2.89%  "

Set-TextValue "D26" "155.60"
Set-TextValue "E26" "  +2.12%  "

Set-TextValue "D27" "19.14"
Set-TextValue "E27" "  +12.82%  "

Set-TextValue "E28" "  +0.08%  "

Set-TextValue "E29" "  +0.93%  "

Set-TextValue "D30" "4.341"
Set-TextValue "E30" "  +2.98%  "

Set-TextValue "D31" "0.09090"
Set-TextValue "E31" "  +3.38%  "

Set-TextValue "E32" "  +0.49%  "

Set-TextValue "D33" "0.05220"
Set-TextValue "E33" "  +2.88%  "

Set-TextValue "D34" "0.7505"
Set-TextValue "E34" "  +3.95%  "

Set-TextValue "D35" "1.109"
Set-TextValue "E35" "  -0.28%  "

Set-TextValue "D37" "0.01839"
Set-TextValue "E37" "  -1.88%  "

Set-TextValue "D38" "2.683"
Set-TextValue "E38" "  -0.15%  "

Set-TextValue "D39" "2.141"
Set-TextValue "E39" "  -0.99%  "

Set-TextValue "D40" "0.9341"
Set-TextValue "E40" "  +0.42%  "

Set-TextValue "D41" "0.4421"
Set-TextValue "E41" "  +4.25%  "

Set-TextValue "D42" "105.42"
Set-TextValue "E42" "  +2.10%  "

Set-TextValue "E43" "  +0.23%  "

Set-TextValue "D44" "5.764"
Set-TextValue "E44" "  -0.48%  "

Set-TextValue "E45" "  +2.79%  "

Set-TextValue "D46" "0.1343"
Set-TextValue "E46" "  +4.94%  "

Set-TextValue "D47" "0.05865"
Set-TextValue "E47" "  +2.90%  "

Set-TextValue "D48" "8.718"
Set-TextValue "E48" "  +5.34%  "

Set-TextValue "D49" "1.425"
Set-TextValue "E49" "  +6.24%  "

Set-TextValue "D50" "0.3941"
Set-TextValue "E50" "  +4.12%  "

Set-TextValue "E51" "  +2.28%  "
